# "readme amended. server cloud config yml centralized"
#
# The "Patch Management tools" row (old row 2: ManageEngine Patch Manager
# Plus + its two hyperlinks) is removed from the "Mine" sheet; every row
# below it shifts up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mine")
$ws.Activate()

# --- 1. stash a clean copy of the plain (non-hyperlink) C/D-column style so
#        we can restore it later — adding a hyperlink auto-applies Excel's
#        built-in "Hyperlink" style, which we don't want here.
$ws.Range("C4").Copy()
$ws.Range("Z1").PasteSpecial(-4122)   # xlPasteFormats

# --- 2. drop the hyperlinks that live in the row about to be removed, then
#        remove the row itself (shifts A3:D8 up into A2:D7).
$ws.Range("C2:D2").Hyperlinks.Delete()
$ws.Rows.Item(2).Delete()

# --- 3. the remaining hyperlinks do not follow their cells on a row
#        delete in this engine, so drop the now-misaligned leftovers and
#        recreate them at their shifted destinations.
$ws.Range("A1:D7").Hyperlinks.Delete()

$links = @(
  @{addr="C1"; url="https://rockylinux.org/download/"; disp="https://rockylinux.org/download/"},
  @{addr="C2"; url="https://ranchermanager.docs.rancher.com/v2.6"; disp="https://ranchermanager.docs.rancher.com/v2.6"},
  @{addr="C3"; url="https://bazel.build/about/intro"; disp="https://bazel.build/about/intro"},
  @{addr="C4"; url="https://buildkite.com/home"; disp="https://buildkite.com/home"},
  @{addr="C5"; url="https://snyk.io/"; disp="https://snyk.io"},
  @{addr="C6"; url="https://spinnaker.io/"; disp="https://spinnaker.io"},
  @{addr="C7"; url="https://sematext.com/pricing/"; disp="https://sematext.com/pricing/"},
  @{addr="D7"; url="https://sematext.com/blog/cloud-monitoring-tools/"; disp="https://sematext.com/blog/cloud-monitoring-tools/"}
)
foreach ($l in $links) {
    $ws.Hyperlinks.Add($ws.Range($l.addr), $l.url, "", "", $l.disp)
}

# --- 4. undo the auto "Hyperlink" styling these Add() calls introduced.
$ws.Range("Z1").Copy()
$ws.Range("C1:C7").PasteSpecial(-4122)
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# --- 5. touching and clearing the final row extends the sheet dimension the
#        same way the source edit did.
$ws.Range("A1048576").Value = "x"
$ws.Range("A1048576").ClearContents()
$ws.Range("A1048576").ClearFormats()

# --- 6. match the recorded selection state.
$ws.Range("A2").Select()
